$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: convert previously-text SmartScore cells to genuine numbers ---
$ws.Cells.Item(8, 9).Value = 0.569
$ws.Cells.Item(8, 12).Value = 0.484
$ws.Cells.Item(8, 15).Value = 0.467
$ws.Cells.Item(8, 18).Value = 0.64
$ws.Cells.Item(8, 21).Value = 0.602
$ws.Cells.Item(8, 24).Value = 0.582
$ws.Cells.Item(8, 27).Value = 0.7
$ws.Cells.Item(8, 30).Value = 0.572
$ws.Cells.Item(8, 33).Value = 0.551

# --- Row 9: new participant (Maribel Badillo) ---
$ws.Cells.Item(9, 1).Value = "Maribel Badillo_20251120_154905"
$ws.Cells.Item(9, 2).Value = "'"
$ws.Cells.Item(9, 3).Value = "Maribel Badillo"
$ws.Cells.Item(9, 4).Value = 21.0
$ws.Cells.Item(9, 5).Value = "Female"
$ws.Cells.Item(9, 6).Value = "2025-11-20 15:49:06"
$ws.Cells.Item(9, 7).Value = @"
{
  "portion": 0.8,
  "diet": 0.14285714285714285,
  "salt": 0.0,
  "fat": 1.0,
  "natural": 0.8,
  "convenience": 0.0,
  "price": 0.4
}
"@
$ws.Cells.Item(9, 8).Value = "Nongshim Neoguri Spicy Seafood"
$ws.Cells.Item(9, 9).Value = "'0.710"
$ws.Cells.Item(9, 10).Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Cells.Item(9, 11).Value = "Maruchan Ramen Sabor Pollo"
$ws.Cells.Item(9, 12).Value = "'0.438"
$ws.Cells.Item(9, 13).Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Cells.Item(9, 14).Value = "Nongshim Shin Ramyun"
$ws.Cells.Item(9, 15).Value = "'0.429"
$ws.Cells.Item(9, 16).Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"
$ws.Cells.Item(9, 17).Value = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item(9, 18).Value = "'0.761"
$ws.Cells.Item(9, 19).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Cells.Item(9, 20).Value = "Annie’s Shells & White Cheddar"
$ws.Cells.Item(9, 21).Value = "'0.706"
$ws.Cells.Item(9, 22).Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Cells.Item(9, 23).Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Cells.Item(9, 24).Value = "'0.517"
$ws.Cells.Item(9, 25).Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Cells.Item(9, 26).Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item(9, 27).Value = "'0.670"
$ws.Cells.Item(9, 28).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Cells.Item(9, 29).Value = "Kitchens of India Variety Pack"
$ws.Cells.Item(9, 30).Value = "'0.582"
$ws.Cells.Item(9, 31).Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"
$ws.Cells.Item(9, 32).Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item(9, 33).Value = "'0.447"
$ws.Cells.Item(9, 34).Value = "Portátil, saludable, fácil, buena textura, sabor suave"

# Restore default (non-custom) row height after writing the multi-line JSON cell,
# matching the source data (which never set an explicit row height).
$ws.Rows.Item(9).EntireRow.AutoFit()
